$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3432038.8
$ws.Range("J9").Value = 3432038.8
$ws.Range("L9").Value = 3432038.8
$ws.Range("N9").Value = -3432376.8
$ws.Range("H12").Value = 324.33334
$ws.Range("I12").Value = 364.2
$ws.Range("J12").Value = 125
$ws.Range("K12").Value = 364.2
$ws.Range("L12").Value = 125
$ws.Range("M12").Value = -194.2
$ws.Range("N12").Value = -465
$ws.Range("H21").Value = 57500836
$ws.Range("I21").Value = 57500836
$ws.Range("K21").Value = 57500836
$ws.Range("M21").Value = -57500368
$ws.Range("H23").Value = 57500836
$ws.Range("I23").Value = 57500836
$ws.Range("K23").Value = 57500836
$ws.Range("M23").Value = -57500602
$ws.Range("H101").Value = 1837.5
$ws.Range("I101").Value = 1798.8889
$ws.Range("K101").Value = 5396.6667
$ws.Range("M101").Value = -3774.6667
$ws.Range("H103").Value = 1761.6428
$ws.Range("I103").Value = 1685.75
$ws.Range("J103").Value = 1862.8334
$ws.Range("K103").Value = 5057.25
$ws.Range("L103").Value = 5588.5002
$ws.Range("M103").Value = -4471.25
$ws.Range("N103").Value = -6760.5002
$ws.Range("H107").Value = 27779408
$ws.Range("I107").Value = 33334068
$ws.Range("J107").Value = 6111.3335
$ws.Range("K107").Value = 33334068
$ws.Range("L107").Value = 6111.3335
$ws.Range("M107").Value = -33332148
$ws.Range("N107").Value = -9951.333500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33497.906
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 223805.4
$ws.Range("I34").Value = 184500
$ws.Range("J34").Value = 250009
$ws.Range("K34").Value = 184500
$ws.Range("L34").Value = 250009
$ws.Range("M34").Value = -184229
$ws.Range("N34").Value = -250551
$ws.Range("H53").Value = 8499.5
$ws.Range("I53").Value = 8499.5
$ws.Range("K53").Value = 8499.5
$ws.Range("M53").Value = -7817.5
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H74").Value = 275210.1
$ws.Range("I74").Value = 419953.97
$ws.Range("J74").Value = 7990.6924
$ws.Range("K74").Value = 419953.97
$ws.Range("L74").Value = 7990.6924
$ws.Range("M74").Value = -419079.97
$ws.Range("N74").Value = -9738.6924
$ws.Range("H77").Value = 275210.1
$ws.Range("I77").Value = 419953.97
$ws.Range("J77").Value = 7990.6924
$ws.Range("K77").Value = 2099769.85
$ws.Range("L77").Value = 39953.462
$ws.Range("M77").Value = -2095401.85
$ws.Range("N77").Value = -48689.462
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H88").Value = 1336.1428
$ws.Range("J88").Value = 1336.1428
$ws.Range("L88").Value = 1336.1428
$ws.Range("N88").Value = -2148.1428
$ws.Range("H91").Value = 1336.1428
$ws.Range("J91").Value = 1336.1428
$ws.Range("L91").Value = 1336.1428
$ws.Range("N91").Value = -4144.1428
$ws.Range("H92").Value = 5025000
$ws.Range("J92").Value = 5025000
$ws.Range("L92").Value = 5025000
$ws.Range("N92").Value = -5029992
$ws.Range("H122").Value = 2706.7693
$ws.Range("I122").Value = 2814
$ws.Range("K122").Value = 8442
$ws.Range("M122").Value = -5992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 754
$ws.Range("I22").Value = 672.8333
$ws.Range("K22").Value = 672.8333
$ws.Range("M22").Value = -322.8333
$ws.Range("H122").Value = 75840.78999999999
$ws.Range("I122").Value = 81528.53999999999
$ws.Range("K122").Value = 244585.62
$ws.Range("M122").Value = -242135.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4279.1665
$ws.Range("I3").Value = 4235
$ws.Range("K3").Value = 12705
$ws.Range("M3").Value = -12593
$ws.Range("H88").Value = 3999.111
$ws.Range("I88").Value = 3999.111
$ws.Range("K88").Value = 11997.333
$ws.Range("M88").Value = -11569.333
$ws.Range("H91").Value = 3999.111
$ws.Range("I91").Value = 3999.111
$ws.Range("K91").Value = 11997.333
$ws.Range("M91").Value = -10515.333
$ws.Range("H136").Value = 857.5714
$ws.Range("I136").Value = 857.5714
$ws.Range("K136").Value = 2572.7142
$ws.Range("M136").Value = 2527.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 3030.5715
$ws.Range("I23").Value = 814
$ws.Range("J23").Value = 3400
$ws.Range("K23").Value = 814
$ws.Range("L23").Value = 3400
$ws.Range("M23").Value = -591
$ws.Range("N23").Value = -3846
$ws.Range("H24").Value = 53027.188
$ws.Range("J24").Value = 36843.5
$ws.Range("L24").Value = 36843.5
$ws.Range("N24").Value = -37189.5
$ws.Range("H97").Value = 1610.25
$ws.Range("I97").Value = 1712.75
$ws.Range("J97").Value = 1405.25
$ws.Range("K97").Value = 1712.75
$ws.Range("L97").Value = 1405.25
$ws.Range("M97").Value = -1216.75
$ws.Range("N97").Value = -2397.25
$ws.Range("H113").Value = 20991.455
$ws.Range("I113").Value = 41532.5
$ws.Range("J113").Value = 9253.714
$ws.Range("K113").Value = 41532.5
$ws.Range("L113").Value = 9253.714
$ws.Range("M113").Value = -39362.5
$ws.Range("N113").Value = -13593.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1997
$ws.Range("I7").Value = 1997
$ws.Range("K7").Value = 1997
$ws.Range("M7").Value = -1885
$ws.Range("H22").Value = 3892.5
$ws.Range("I22").Value = 1729.1428
$ws.Range("J22").Value = 4783.294
$ws.Range("K22").Value = 1729.1428
$ws.Range("L22").Value = 4783.294
$ws.Range("M22").Value = -1434.1428
$ws.Range("N22").Value = -5373.294
$ws.Range("H27").Value = 3892.5
$ws.Range("I27").Value = 1729.1428
$ws.Range("J27").Value = 4783.294
$ws.Range("K27").Value = 1729.1428
$ws.Range("L27").Value = 4783.294
$ws.Range("M27").Value = -1622.1428
$ws.Range("N27").Value = -4997.294
$ws.Range("H55").Value = 620.7857
$ws.Range("I55").Value = 878.8570999999999
$ws.Range("J55").Value = 362.7143
$ws.Range("K55").Value = 878.8570999999999
$ws.Range("L55").Value = 362.7143
$ws.Range("M55").Value = -705.8570999999999
$ws.Range("N55").Value = -708.7143
$ws.Range("H69").Value = 44888
$ws.Range("J69").Value = 44888
$ws.Range("L69").Value = 44888
$ws.Range("N69").Value = -46510
$ws.Range("H72").Value = 44888
$ws.Range("J72").Value = 44888
$ws.Range("L72").Value = 134664
$ws.Range("N72").Value = -142776
$ws.Range("H82").Value = 2280.7
$ws.Range("J82").Value = 2403.1428
$ws.Range("L82").Value = 2403.1428
$ws.Range("N82").Value = -3125.1428
$ws.Range("H85").Value = 2280.7
$ws.Range("J85").Value = 2403.1428
$ws.Range("L85").Value = 2403.1428
$ws.Range("N85").Value = -4899.1428
$ws.Range("H126").Value = 1997
$ws.Range("I126").Value = 1997
$ws.Range("K126").Value = 5991
$ws.Range("M126").Value = -3521
$ws.Range("H132").Value = 4340.591
$ws.Range("I132").Value = 3074.6875
$ws.Range("J132").Value = 7716.3335
$ws.Range("K132").Value = 9224.0625
$ws.Range("L132").Value = 23149.0005
$ws.Range("M132").Value = -6694.0625
$ws.Range("N132").Value = -28209.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H96").Value = 1078.3334
$ws.Range("J96").Value = 1350
$ws.Range("L96").Value = 1350
$ws.Range("N96").Value = -4096
$ws.Range("H107").Value = 2488.6365
$ws.Range("I107").Value = 1943.5
$ws.Range("K107").Value = 5830.5
$ws.Range("M107").Value = -3910.5
$ws.Range("H122").Value = 1695
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1695
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 5085
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9985
$ws.Range("H126").Value = 2292.7932
$ws.Range("I126").Value = 1433.72
$ws.Range("J126").Value = 7662
$ws.Range("K126").Value = 4301.16
$ws.Range("L126").Value = 22986
$ws.Range("M126").Value = -1831.16
$ws.Range("N126").Value = -27926
$ws.Range("H132").Value = 3879.4443
$ws.Range("I132").Value = 1756.1
$ws.Range("J132").Value = 6533.625
$ws.Range("K132").Value = 5268.299999999999
$ws.Range("L132").Value = 19600.875
$ws.Range("M132").Value = -2738.299999999999
$ws.Range("N132").Value = -24660.875
